# "SCADA: Inicio y avance de programacion"
#
# Updates the PWM sizing worksheet (Hoja1):
#   - Top "datos conocidos" block (row 5): new Pre-escaler / F(PWM) / Duty inputs.
#   - TIMER2 block (rows 12-14): new Fosc / Pre-escaler(D12) / bits(E14) inputs.
#   - New F(PWM) / period read-outs (E17/E18) added under the TIMER0 block.
#   - TIMER0 block (rows 19): new Fosc / Pre-escaler / F(PWM) inputs.
#   - Row 16 grows a bit taller; row 15 loses its explicit custom height.
#   - Selection cursor moves to C19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- "SE CONOCEN ESTOS DATOS" block (row 5) ---
# B5 (Fosc) is left as-is.
$ws.Range("C5").Value = 1          # Pre
$ws.Range("D5").Value = 200000     # Fpwm
$ws.Range("E5").Value = 1          # Duty

# --- TIMER2 block (rows 12-14) ---
$ws.Range("B12").Value = 20        # Fosc (Mhz)
$ws.Range("D12").Value = 32        # Pre-escaler Timer 2
$ws.Range("E12").Value = 1000      # F(PWM) Hz
$ws.Range("E14").Value = 8         # Bits de Timer

# --- New read-outs under the TIMER0 header (rows 17-18) ---
$ws.Range("E17").Formula = "=1/E18"
$ws.Range("E18").Formula = "=4/(B12*1000000)*_xlfn.BITLSHIFT(1,E14)*D12"

# --- TIMER0 block (row 19) ---
$ws.Range("B19").Value = 20        # Fosc (Mhz)
$ws.Range("C19").Value = 4         # Pre-escaler Timer 0
$ws.Range("D19").Value = 500000    # Frecuencia (Hz)

# --- Row height tweaks ---
$ws.Rows.Item(16).RowHeight = 21

# --- Window / selection state ---
try {
    $win = $wb.Windows.Item(1)
    $win.WindowState = -4140  # xlMinimized
    $win.Left = 30
    $win.Top = 750
    $win.Width = 28770
    $win.Height = 15450
} catch {}

$ws.Range("C19").Select()

$wb.Save()
